# Edit LOM3230.xlsx: remove the stray "Docentes responsáveis:" value row
# (old row 13, which only held the teacher name in columns B/C) by deleting
# it -- this shifts every following row up by one and keeps the per-row
# heights/styles intact. Then fix up the handful of cells whose text content
# changed as a result of the underlying data re-shuffle described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 13 (B13/C13 = "5840726 - Cristina Bormio Nunes").
# Everything below shifts up by one row, preserving row heights/styles.
$ws.Rows("13").Delete()

# Row 10 (Objetivos:) now shows the teacher id/name instead of the old text.
$ws.Range("B10").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C10").Value = "5840726 - Cristina Bormio Nunes"

# Row 13 (Programa resumido:, was old row 14) now reads "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:, was old row 16) now reads "01/01/2012". Entering this
# literal via .Value would be auto-parsed as a date serial number by Excel,
# so instead put a text formula on the cell, then paste-special just the
# (already textual) value back over itself -- this keeps it a plain string
# without touching the cell's existing number format/style.
$b15 = $ws.Range("B15")
$b15.Formula = "=""01/01/2012"""
$b15.Copy()
$b15.PasteSpecial(-4163)

$c15 = $ws.Range("C15")
$c15.Formula = "=""01/01/2012"""
$c15.Copy()
$c15.PasteSpecial(-4163)

# Row 18 (Método:, was old row 19) now shows the teacher id/name.
$ws.Range("B18").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C18").Value = "5840726 - Cristina Bormio Nunes"

# Row 19 (Critério:, was old row 20) now holds the lab-experiments text.
$ws.Range("B19").Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."
$ws.Range("C19").Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."

# Row 20 (Norma de recuperação:, was old row 21) now holds the grading formula text.
$ws.Range("B20").Value = "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3"
$ws.Range("C20").Value = "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3"

# Row 21 (Bibliografia:, was old row 22) now holds the recovery-exam text.
$ws.Range("B21").Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
